# Atualização de casos/óbitos até 22/04/22
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(44665, 0, 326862, 6338, 28, 0),
    @(44666, 0, 326879, 6339, 17, 1),
    @(44667, 0, 326894, 6339, 15, 0),
    @(44668, 0, 326902, 6339, 8, 0),
    @(44669, 0, 326909, 6339, 7, 0),
    @(44670, 0, 326928, 6340, 19, 1),
    @(44671, 0, 326949, 6341, 21, 1),
    @(44672, 0, 326970, 6341, 21, 0),
    @(44673, 0, 326989, 6342, 19, 1)
)

$startRow = 19
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}

# Scroll the view so row 5 is at the top (best-effort; mirrors the
# author's saved scroll position) without disturbing the final selection.
try {
    $win = $wb.Windows.Item(1)
    $win.ScrollRow = 5
    $win.ScrollColumn = 1
} catch {
}

$ws.Range("E21").Select()
